$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B updates (rows 2-15) ---
$ws.Range("B2").Value = "NSE:ADL"
$ws.Range("B3").Value = "NSE:ASTRAMICRO"
$ws.Range("B4").Value = "NSE:AVTNPL"
$ws.Range("B5").Value = "NSE:BAJAJ-AUTO"
$ws.Range("B6").Value = "NSE:BINANIIND"
$ws.Range("B7").Value = "NSE:BRITANNIA"
$ws.Range("B8").Value = "NSE:CAMPUS"
$ws.Range("B9").Value = "NSE:CGPOWER"
$ws.Range("B10").Value = "NSE:CONSUMBEES"
$ws.Range("B11").Value = "NSE:GANDHITUBE"
$ws.Range("B12").Value = "NSE:GOKULAGRO"
$ws.Range("B13").Value = "NSE:INDIGO"
$ws.Range("B14").Value = "NSE:INDORAMA"
$ws.Range("B15").Value = "NSE:KAYNES"

# --- Column C updates (rows 3-15) ---
$ws.Range("C3").Value = "NSE:AGI"
$ws.Range("C5").Value = "NSE:ALKYLAMINE"
$ws.Range("C6").Value = "NSE:ASIANENE"
$ws.Range("C7").Value = "NSE:GINNIFILA"
$ws.Range("C8").Value = "NSE:HINDWAREAP"
$ws.Range("C9").Value = "NSE:KERNEX"
$ws.Range("C10").Value = "NSE:LAL"
$ws.Range("C11").Value = "NSE:OAL"
$ws.Range("C12").ClearContents()
$ws.Range("C13").ClearContents()
$ws.Range("C14").ClearContents()
$ws.Range("C15").ClearContents()

# --- Column E updates (rows 2-3 cleared) ---
$ws.Range("E2").ClearContents()
$ws.Range("E3").ClearContents()

# --- Column F updates (row 3 now has value) ---
$ws.Range("F3").Value = "NSE:ITC"

# --- New rows 16-22 ---
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "NSE:KSL"

$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "NSE:LAOPALA"

$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "NSE:MAPMYINDIA"

$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "NSE:MEDICAMEQ"

$ws.Range("A20").Value = 18
$ws.Range("B20").Value = "NSE:MOHEALTH"

$ws.Range("A21").Value = 19
$ws.Range("B21").Value = "NSE:PAGEIND"

$ws.Range("A22").Value = 20
$ws.Range("B22").Value = "NSE:PREMIERPOL"

# Apply same style as other A column cells (style index 1, bold/centered/bordered)
# to the new A column cells by copying formatting from A15
$ws.Range("A15").Copy()
$ws.Range("A16:A22").PasteSpecial(-4122)
$excel.CutCopyMode = $false
